$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Preserve the "Total" row's bold/border formatting (currently on row 32)
#        by relocating it to row 55 *before* row 32 gets overwritten below. ---
$ws.Range("B32:C32").Copy()
$ws.Range("B55:C55").PasteSpecial(-4122)  # xlPasteFormats

# --- 2. Extend formatting (styles s="3"/"4"/"8" for date/activity/hours columns)
#        down through row 54, matching the blank placeholder rows 28-31 that
#        already carry that formatting in the source file. ---
$ws.Range("A31:C31").Copy()
$ws.Range("A32:C54").PasteSpecial(-4122)  # xlPasteFormats

# --- 3. Fill in the new journal entries (rows 28-37) ---

# Row 28
$ws.Range("A28").Value = 43199
$ws.Range("B28").Value = "Création de la fenêtre de création d'un compte bancaire, Problème d'interaction avec la fenêtre principale. Je n'arrive pas a faire communiquer les deux controlleur (Controlleur_BankAccount et Controlleur_createBankAccount)"
$ws.Range("C28").Value = 1

# Row 29
$ws.Range("A29").Value = 43203
$ws.Range("B29").Value = "Résolution du problème javaFX et maven"
$ws.Range("C29").Value = 1.5

# Row 30
$ws.Range("A30").Value = 43203
$ws.Range("B30").Value = "Intégration du javaFX dans le projet maven, problème : il n'arrive pas a load les fichier fxml. Recherche du problème"
$ws.Range("C30").Value = 1.5

# Row 31
$ws.Range("A31").Value = 43205
$ws.Range("B31").Value = "création du ppt"
$ws.Range("C31").Value = 0.5

# Row 32
$ws.Range("A32").Value = 43206
$ws.Range("B32").Value = "présentation intermédiaire de notre projet"
$ws.Range("C32").Value = 0.25

# Row 33
$ws.Range("A33").Value = 43206
$ws.Range("B33").Value = "Ajout des fonction de la BLL dans le controlleur du loginRegister"
$ws.Range("C33").Value = 0.75

# Row 34
$ws.Range("A34").Value = 43207
$ws.Range("B34").Value = "création d'un compte utilisateur terminé"
$ws.Range("C34").Value = 0.25

# Row 35
$ws.Range("A35").Value = 43213
$ws.Range("B35").Value = "Mise au point de ce qu'il reste a faire sur l'interface graphique"
$ws.Range("C35").Value = 1.5

# Row 36
$ws.Range("A36").Value = 43218
$ws.Range("B36").Value = "Modification du rapport de la partie controlleur du compte bancaire"
$ws.Range("C36").Value = 0.25

# Row 37
$ws.Range("A37").Value = 43218
$ws.Range("B37").Value = "Ajout de l'évenement sur le compte bancaire"
$ws.Range("C37").Value = 0.5

# --- 4. Row heights for wrapped multi-line entries ---
$ws.Rows.Item(28).RowHeight = 45
$ws.Rows.Item(30).RowHeight = 75
$ws.Rows.Item(33).RowHeight = 30
$ws.Rows.Item(35).RowHeight = 30
$ws.Rows.Item(36).RowHeight = 30

# --- 5. Write the "Total" row content/formula in its new location (row 55) ---
$ws.Range("B55").Value = "Total"
$ws.Range("C55").Formula = "=SUM(C5:C54)"

# --- 6. Update sheet view: scroll position & selection ---
$ws.Range("C38").Select()
